# Add the new "Sprint 2" log entry row (row 14) describing the admin
# functionality work, matching the existing table's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Admin page, connection to db, adding and removing brokers from website. Working on editing broker"
$ws.Range("B14").Value = "Contacted kumai a lot "
$ws.Range("C14").Value = "Gabriel and Kumai"
$ws.Range("D14").Value = "9h"

# Copy the formatting used by the other detail rows (e.g. row 4) onto the
# new row so it matches the table's look (fill, border, centered/wrapped).
$ws.Range("A4:D4").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)

# Reflect where the user ended up after adding the row: scrolled down and
# sitting on the first empty row below the new entry.
$ws.Range("A15").Select()
